$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Add single-line borders (top/left/bottom/right, sz=4, color=auto)
#    to every cell of the table (w:tcBorders inside each w:tcPr).
# ---------------------------------------------------------------------
$t = $d.Tables(1)
$sides = -1, -2, -3, -4
foreach ($row in $t.Rows) {
    foreach ($cell in $row.Cells) {
        $cell.Borders.DistanceFromTop = 0
        $cell.Borders.DistanceFromLeft = 0
        $cell.Borders.DistanceFromBottom = 0
        $cell.Borders.DistanceFromRight = 0
        foreach ($s in $sides) {
            $b = $cell.Borders.Item($s)
            $b.LineStyle = 1
            $b.LineWidth = 2
            $b.ColorIndex = 0
        }
    }
}

# ---------------------------------------------------------------------
# 2) Relocate the reserved "_GoBack" bookmark: it currently sits in the
#    empty paragraph right after the table; it must instead wrap the
#    start of the paragraph holding the "m:usercontent zone1" field
#    (the paragraph right before the table). Re-adding a bookmark named
#    "_GoBack" moves it (Word keeps a single "last edit" bookmark).
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$fieldParaStart = -1
$bookmarkParaStart = -1
$bookmarkParaEnd = -1

$i = 0
foreach ($p in $paras) {
    $i = $i + 1
    if ($i -eq 2) {
        $fieldParaStart = $p.Range.Start
    }
    if ($i -eq 23) {
        $bookmarkParaStart = $p.Range.Start
        $bookmarkParaEnd = $p.Range.End
    }
}

$insertionPoint = $d.Range($fieldParaStart, $fieldParaStart)
$d.Bookmarks.Add("_GoBack", $insertionPoint)

# ---------------------------------------------------------------------
# 3) Remove the now bookmark-less empty paragraph that used to follow
#    the table (it contained only the bookmark and an empty run).
# ---------------------------------------------------------------------
$deadPara = $d.Range($bookmarkParaStart, $bookmarkParaEnd)
$deadPara.Delete()
